# fix issue with DS 6 results
#
# Re-run (corrected) Kendall's tau / Spearman's rho correlation figures for the
# dataset_id == 6 rows (rows 10-12, i.e. "DS 6") on every results sheet.
# Columns I/J/K/L hold kendalls_tau, kendalls_p_value, spearmans_rho and
# spearmans_p_value respectively.
#
# A couple of sheets also had column 9/11 slightly mis-sized (21.71 chars
# instead of 20.71 chars, i.e. one character too wide) - align them with the
# neighbouring columns while we're in here.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# all_tools
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("all_tools")
$ws.Columns.Item(11).ColumnWidth = 19.83

$ws.Range("I10").Value = -0.004414751593059719
$ws.Range("J10").Value = 0.9661430518463994
$ws.Range("K10").Value = 0.01233315619210278
$ws.Range("L10").Value = 0.932256616388448

$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = 0.01007744778718968
$ws.Range("L11").Value = 0.9446252971798705

$ws.Range("I12").Value = 0.04719010357797921
$ws.Range("J12").Value = 0.632959082865391
$ws.Range("K12").Value = 0.07620564634846749
$ws.Range("L12").Value = 0.5988937623118764

# ---------------------------------------------------------------------------
# checker_framework
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("checker_framework")

$ws.Range("I10").Value = -0.2741939043543897
$ws.Range("J10").Value = 0.01523742906785629
$ws.Range("K10").Value = -0.3438060590640694
$ws.Range("L10").Value = 0.0145005762354219

$ws.Range("I11").Value = -0.2885159533670518
$ws.Range("J11").Value = 0.007534224862278689
$ws.Range("K11").Value = -0.3916527812158478
$ws.Range("L11").Value = 0.004913052567611288

$ws.Range("I12").Value = 0.1527404930393114
$ws.Range("J12").Value = 0.1545006958890939
$ws.Range("K12").Value = 0.1874148805812322
$ws.Range("L12").Value = 0.192476549579598

# ---------------------------------------------------------------------------
# typestate_checker
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("typestate_checker")
$ws.Columns.Item(9).ColumnWidth = 19.83
$ws.Columns.Item(11).ColumnWidth = 19.83

$ws.Range("I10").Value = 0.03324642499485375
$ws.Range("J10").Value = 0.7524695495544889
$ws.Range("K10").Value = 0.05369692827151148
$ws.Range("L10").Value = 0.7111154281850389

$ws.Range("I11").Value = 0.02630142666529713
$ws.Range("J11").Value = 0.7941029996140164
$ws.Range("K11").Value = 0.0345247475921612
$ws.Range("L11").Value = 0.8118632018166435

$ws.Range("I12").Value = 0.02696081483975262
$ws.Range("J12").Value = 0.7877495099042957
$ws.Range("K12").Value = 0.03904490288887924
$ws.Range("L12").Value = 0.7877690834629696

# ---------------------------------------------------------------------------
# infer
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("infer")
$ws.Columns.Item(11).ColumnWidth = 19.83

$ws.Range("I10").Value = -0.1454025530693833
$ws.Range("J10").Value = 0.2372373518450496
$ws.Range("K10").Value = -0.17271903862684
$ws.Range("L10").Value = 0.2303502122764337

$ws.Range("I11").Value = -0.1395616700784287
$ws.Range("J11").Value = 0.2348980869048207
$ws.Range("K11").Value = -0.1674579385094694
$ws.Range("L11").Value = 0.2450782275649824

$ws.Range("I12").Value = 0.02140819589682411
$ws.Range("J12").Value = 0.8544862615484419
$ws.Range("K12").Value = 0.02708713119452734
$ws.Range("L12").Value = 0.8518765230635053

# ---------------------------------------------------------------------------
# openjml
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("openjml")

$ws.Range("I10").Value = -0.1679543721361666
$ws.Range("J10").Value = 0.1236240407674998
$ws.Range("K10").Value = -0.1891736836680379
$ws.Range("L10").Value = 0.1929773116581186

$ws.Range("I11").Value = -0.1257020377320922
$ws.Range("J11").Value = 0.2278618094667795
$ws.Range("K11").Value = -0.1630362203611797
$ws.Range("L11").Value = 0.2630118256866814

$ws.Range("I12").Value = 0.09860866138702759
$ws.Range("J12").Value = 0.3411530572805948
$ws.Range("K12").Value = 0.1348255148415411
$ws.Range("L12").Value = 0.3556729989431554
